$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("r2 (2Ysum)")
$dst = $wb.Worksheets.Item("r2 (3Ysum)")

$src.Range("S3:T57").Copy()
$dst.Range("T2:U56").PasteSpecial(-4122)
$src.Range("S3:T57").Copy()
$dst.Range("Y2:Z56").PasteSpecial(-4122)

Write-Host "T2 style:" $dst.Range("T2").Interior.Color
Write-Host "U2 style:" $dst.Range("U2").Interior.Color
Write-Host "T4 style:" $dst.Range("T4").Interior.Color  $dst.Range("T4").Font.Bold
Write-Host "U4 style:" $dst.Range("U4").Interior.Color $dst.Range("U4").NumberFormat
Write-Host "T54 style:" $dst.Range("T54").Interior.Color
Write-Host "T55 style:" $dst.Range("T55").Interior.Color $dst.Range("T55").Font.Bold
